$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.06440233333333333
$ws.Range("H2").Value = 0.193207
$ws.Range("I2").Value = 0.03647206354366116
$ws.Range("J2").Value = 0.03647206354366116
$ws.Range("M2").Value = 1.443038
$ws.Range("N2").Value = 4.329114
$ws.Range("O2").Value = 0.0289666880885598
$ws.Range("P2").Value = 0.0289666880885598
$ws.Range("Q2").Value = 0.09293501428866664
$ws.Range("R2").Value = 0.8364151285979999
$ws.Range("S2").Value = 0.001056474888615366
$ws.Range("T2").Value = 0.001056474888615366
# Row 3
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.06440233333333333
$ws.Range("H3").Value = 0.193207
$ws.Range("I3").Value = 0.03647206354366116
$ws.Range("J3").Value = 0.03647206354366116
$ws.Range("N3").Value = 87.610543
$ws.Range("O3").Value = 0.5862140087672342
$ws.Range("P3").Value = 0.5862140087672342
$ws.Range("Q3").Value = 1.880774464600111
$ws.Range("R3").Value = 16.926970181401
$ws.Range("S3").Value = 0.0213804345779429
$ws.Range("T3").Value = 0.0213804345779429
# Row 4
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.06440233333333333
$ws.Range("H4").Value = 0.193207
$ws.Range("I4").Value = 0.03647206354366116
$ws.Range("J4").Value = 0.03647206354366116
$ws.Range("M4").Value = 19.170603
$ws.Range("N4").Value = 57.511809
$ws.Range("O4").Value = 0.384819303144206
$ws.Range("P4").Value = 0.384819303144206
$ws.Range("Q4").Value = 1.234631564607
$ws.Range("R4").Value = 11.111684081463
$ws.Range("S4").Value = 0.01403515407710289
$ws.Range("T4").Value = 0.01403515407710289
# Row 5
$ws.Range("I5").Value = 0.8194013021867156
$ws.Range("J5").Value = 0.8194013021867155
$ws.Range("M5").Value = 1.443038
$ws.Range("N5").Value = 4.329114
$ws.Range("O5").Value = 0.0289666880885598
$ws.Range("P5").Value = 0.0289666880885598
$ws.Range("Q5").Value = 2.087928796124
$ws.Range("R5").Value = 18.791359165116
$ws.Range("S5").Value = 0.02373534193980232
$ws.Range("T5").Value = 0.02373534193980232
# Row 6
$ws.Range("I6").Value = 0.8194013021867156
$ws.Range("J6").Value = 0.8194013021867155
$ws.Range("N6").Value = 87.610543
$ws.Range("O6").Value = 0.5862140087672342
$ws.Range("P6").Value = 0.5862140087672342
$ws.Range("Q6").Value = 42.25450648187134
$ws.Range("R6").Value = 380.290558336842
$ws.Range("S6").Value = 0.4803445221439664
$ws.Range("T6").Value = 0.4803445221439663
# Row 7
$ws.Range("I7").Value = 0.8194013021867156
$ws.Range("J7").Value = 0.8194013021867155
$ws.Range("M7").Value = 19.170603
$ws.Range("N7").Value = 57.511809
$ws.Range("O7").Value = 0.384819303144206
$ws.Range("P7").Value = 0.384819303144206
$ws.Range("Q7").Value = 27.737907139494
$ws.Range("R7").Value = 249.641164255446
$ws.Range("S7").Value = 0.3153214381029469
$ws.Range("T7").Value = 0.3153214381029468
# Row 8
$ws.Range("E8").Value = 1
$ws.Range("F8").Value = 0.3333333333333333
$ws.Range("G8").Value = 0.2544986666666667
$ws.Range("H8").Value = 0.763496
$ws.Range("I8").Value = 0.1441266342696234
$ws.Range("J8").Value = 0.1441266342696234
$ws.Range("M8").Value = 1.443038
$ws.Range("N8").Value = 4.329114
$ws.Range("O8").Value = 0.0289666880885598
$ws.Range("P8").Value = 0.0289666880885598
$ws.Range("Q8").Value = 0.3672512469493333
$ws.Range("R8").Value = 3.305261222543999
$ws.Range("S8").Value = 0.004174871260142114
$ws.Range("T8").Value = 0.004174871260142113
# Row 9
$ws.Range("E9").Value = 1
$ws.Range("F9").Value = 0.3333333333333333
$ws.Range("G9").Value = 0.2544986666666667
$ws.Range("H9").Value = 0.763496
$ws.Range("I9").Value = 0.1441266342696234
$ws.Range("J9").Value = 0.1441266342696234
$ws.Range("N9").Value = 87.610543
$ws.Range("O9").Value = 0.5862140087672342
$ws.Range("P9").Value = 0.5862140087672342
$ws.Range("Q9").Value = 7.432255459814222
$ws.Range("R9").Value = 66.890299138328
$ws.Range("S9").Value = 0.08448905204532493
$ws.Range("T9").Value = 0.08448905204532493
# Row 10
$ws.Range("E10").Value = 1
$ws.Range("F10").Value = 0.3333333333333333
$ws.Range("G10").Value = 0.2544986666666667
$ws.Range("H10").Value = 0.763496
$ws.Range("I10").Value = 0.1441266342696234
$ws.Range("J10").Value = 0.1441266342696234
$ws.Range("M10").Value = 19.170603
$ws.Range("N10").Value = 57.511809
$ws.Range("O10").Value = 0.384819303144206
$ws.Range("P10").Value = 0.384819303144206
$ws.Range("Q10").Value = 4.878892902695999
$ws.Range("R10").Value = 43.910036124264
$ws.Range("S10").Value = 0.0554627109641563
$ws.Range("T10").Value = 0.05546271096415629
